$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions Tracker")

# --- Row 33: partially filled in (Topic continuation row) ---
$ws.Range("B33").Value = "Total number of Occurance"
$ws.Range("H33").Value = "Same as above question. Just substract last minus first."

# --- Row 34: new question entry, fully filled in ---
$ws.Range("B34").Value = "852. Peak Index in a Mountain Array"

# C34 needs the same "Hyperlink" look-and-feel (style) as the other Link
# cells (e.g. C32) plus an actual external hyperlink relationship.
# Copying format from an existing hyperlink cell keeps us on the already
# existing "Hyperlink" cell style instead of minting a brand new one.
$ws.Range("C32").Copy($ws.Range("C34"))
$ws.Range("C34").Value = "LeetCode"
$ws.Hyperlinks.Add($ws.Range("C34"), "https://leetcode.com/problems/peak-index-in-a-mountain-array/") | Out-Null
# Re-apply the shared Hyperlink formatting (Hyperlinks.Add forces its own
# style) so the cell keeps using the existing Hyperlink cell style.
$ws.Range("C32").Copy($ws.Range("C34"))
$ws.Range("C34").Value = "LeetCode"

$ws.Range("D34").Value = 1
$ws.Range("E34").Value = "14/11/2022"
$ws.Range("F34").Value = 3
$ws.Range("G34").Value = "Yes"
$ws.Range("H34").Value = "Good Question"

# --- Update the saved selection / active cell on the sheet ---
$ws.Range("B35").Select() | Out-Null
